{"js": "// Update the worksheet date and the 25 three-digit-by-one-digit\n// multiplication prompts to the new values from the commit.\nconst replacements = [\n  [\"2024-05-11 Saturday\", \"2024-05-12 Sunday\"],\n  [\"592\u00d76=\", \"589\u00d79=\"],\n  [\"175\u00d72=\", \"354\u00d74=\"],\n  [\"313\u00d75=\", \"947\u00d74=\"],\n  [\"342\u00d76=\", \"816\u00d74=\"],\n  [\"579\u00d74=\", \"719\u00d76=\"],\n  [\"558\u00d77=\", \"257\u00d72=\"],\n  [\"998\u00d72=\", \"391\u00d78=\"],\n  [\"667\u00d72=\", \"158\u00d77=\"],\n  [\"877\u00d78=\", \"475\u00d78=\"],\n  [\"716\u00d76=\", \"781\u00d79=\"],\n  [\"536\u00d74=\", \"720\u00d78=\"],\n  [\"285\u00d73=\", \"362\u00d76=\"],\n  [\"997\u00d73=\", \"794\u00d78=\"],\n  [\"453\u00d73=\", \"231\u00d76=\"],\n  [\"663\u00d72=\", \"286\u00d73=\"],\n  [\"673\u00d74=\", \"726\u00d78=\"],\n  [\"323\u00d79=\", \"575\u00d78=\"],\n  [\"902\u00d73=\", \"909\u00d73=\"],\n  [\"455\u00d77=\", \"745\u00d74=\"],\n  [\"443\u00d78=\", \"564\u00d78=\"],\n  [\"692\u00d73=\", \"835\u00d72=\"],\n  [\"304\u00d77=\", \"524\u00d75=\"],\n  [\"821\u00d79=\", \"452\u00d78=\"],\n  [\"732\u00d73=\", \"817\u00d72=\"],\n  [\"961\u00d79=\", \"466\u00d78=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 three-digit-by-one-digit\n# multiplication prompts to the new values from the commit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2024-05-11 Saturday', '2024-05-12 Sunday'),\n    @('592\u00d76=', '589\u00d79='),\n    @('175\u00d72=', '354\u00d74='),\n    @('313\u00d75=', '947\u00d74='),\n    @('342\u00d76=', '816\u00d74='),\n    @('579\u00d74=', '719\u00d76='),\n    @('558\u00d77=', '257\u00d72='),\n    @('998\u00d72=', '391\u00d78='),\n    @('667\u00d72=', '158\u00d77='),\n    @('877\u00d78=', '475\u00d78='),\n    @('716\u00d76=', '781\u00d79='),\n    @('536\u00d74=', '720\u00d78='),\n    @('285\u00d73=', '362\u00d76='),\n    @('997\u00d73=', '794\u00d78='),\n    @('453\u00d73=', '231\u00d76='),\n    @('663\u00d72=', '286\u00d73='),\n    @('673\u00d74=', '726\u00d78='),\n    @('323\u00d79=', '575\u00d78='),\n    @('902\u00d73=', '909\u00d73='),\n    @('455\u00d77=', '745\u00d74='),\n    @('443\u00d78=', '564\u00d78='),\n    @('692\u00d73=', '835\u00d72='),\n    @('304\u00d77=', '524\u00d75='),\n    @('821\u00d79=', '452\u00d78='),\n    @('732\u00d73=', '817\u00d72='),\n    @('961\u00d79=', '466\u00d78='),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair[1]\n    $null = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
